$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-breaking space (U+00A0), as used between "$" and the price in the
# existing rows (e.g. row 4/5: "... $<nbsp> 16.02 ...").
$nbsp = [char]0x00A0

$ws.Range("A6").Value = "r1ange white angel hair 6/16-oz `$$nbsp 16.02 1"
$ws.Range("B6").Value = "1 01ANGE"
$ws.Range("C6").Value = ""

$ws.Range("A7").Value = "r1ling white linguini 6/16-oz `$$nbsp 16.02 2"
$ws.Range("B7").Value = "2 01LING"
$ws.Range("C7").Value = ""
